$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-13 Tuesday" "2024-02-14 Wednesday"

Replace-Text "160÷4=" "979÷5="
Replace-Text "629÷8=" "256÷6="
Replace-Text "991÷6=" "470÷4="
Replace-Text "890÷8=" "473÷4="
Replace-Text "491÷4=" "301÷7="

Replace-Text "258÷4=" "795÷3="
Replace-Text "334÷4=" "802÷4="
Replace-Text "908÷5=" "810÷9="
Replace-Text "652÷2=" "286÷7="
Replace-Text "337÷2=" "715÷8="

Replace-Text "528÷9=" "512÷3="
Replace-Text "483÷8=" "345÷9="
Replace-Text "545÷3=" "516÷2="
Replace-Text "969÷6=" "218÷6="
Replace-Text "231÷8=" "164÷7="

Replace-Text "378÷2=" "226÷7="
Replace-Text "986÷2=" "449÷6="
Replace-Text "880÷8=" "314÷8="
Replace-Text "256÷3=" "825÷5="
Replace-Text "404÷7=" "236÷7="

Replace-Text "578÷2=" "284÷9="
Replace-Text "330÷8=" "509÷9="
Replace-Text "779÷8=" "967÷3="
Replace-Text "731÷8=" "627÷5="
Replace-Text "298÷8=" "345÷4="
